# Apply cell value updates per the diff (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.526.72"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.903.06"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'246.54"
$ws.Range("E5").Value = "  +5.51%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'42.21"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "'0.0703"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "2.177.68"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").Value = "'12.39"
$ws.Range("E13").Value = "  +7.51%  "
$ws.Range("D14").Value = "1.905.20"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "'0.690"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").Value = "'4.86"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").Value = "35.483.86"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "'71.86"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "0.0₃0812"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "'243.39"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'4.88"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "'2.29"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +32.92%  "
$ws.Range("D26").Value = "'171.75"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'8.59"
$ws.Range("E27").Value = "  +8.11%  "
$ws.Range("D28").Value = "'17.96"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "'0.977"
$ws.Range("E30").Value = "  +27.62%  "
$ws.Range("D31").Value = "'4.09"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").Value = "'0.0565"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "'4.14"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("E35").Value = "  +5.98%  "
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  +5.64%  "
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("B40").Value = "MultiversX"
$ws.Range("C40").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D40").Value = "'51.73"
$ws.Range("E40").Value = "  +50.00%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'91.01"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "1.354.14"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'15.49"
$ws.Range("E43").Value = "  +5.57%  "
$ws.Range("D44").Value = "'0.0594"
$ws.Range("E44").Value = "  +11.41%  "
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "'12.60"
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").Value = "'6.66"
$ws.Range("E49").Value = "  +4.63%  "
$ws.Range("D50").Value = "2.087.18"
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("E51").Value = "  +2.18%  "
